$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "init" parameter column (column I), mirroring column H's layout
$ws.Range("I1").Value = "init"
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1050
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 40
$ws.Range("I8").Value = 50
$ws.Range("I9").Value = 200
$ws.Range("I10").Value = 0
$ws.Range("I12").Value = 0.5
$ws.Range("I13").Value = 0
$ws.Range("I14").Value = 1

# Update selection to match the new active cell
$ws.Range("I10").Select()
